# Update "想去人数" (F column) counts across sheets to reflect refreshed
# gh-pages data output.

$wb = $excel.ActiveWorkbook

$ws1 = $wb.Worksheets.Item("展览")
$ws1Vals = @{
    5  = 1306
    6  = 253
    7  = 389
    8  = 8370
    10 = 10301
    21 = 70
    23 = 400
    25 = 1756
    26 = 58
    27 = 518
    28 = 333
    30 = 54
    33 = 1093
    42 = 498
    48 = 63
    49 = 66
}
foreach ($row in $ws1Vals.Keys) {
    $ws1.Range("F$row").Value = $ws1Vals[$row]
}

$ws2 = $wb.Worksheets.Item("演出")
$ws2Vals = @{
    19 = 373
}
foreach ($row in $ws2Vals.Keys) {
    $ws2.Range("F$row").Value = $ws2Vals[$row]
}

$ws3 = $wb.Worksheets.Item("本地生活")
$ws3Vals = @{
    3 = 2782
    4 = 337
    5 = 202
}
foreach ($row in $ws3Vals.Keys) {
    $ws3.Range("F$row").Value = $ws3Vals[$row]
}

$ws4 = $wb.Worksheets.Item("全部类型")
$ws4Vals = @{
    5  = 202
    7  = 1306
    9  = 253
    12 = 389
    13 = 8370
    15 = 10301
    21 = 70
    23 = 1756
    24 = 518
    25 = 333
    30 = 1093
    40 = 498
    47 = 373
}
foreach ($row in $ws4Vals.Keys) {
    $ws4.Range("F$row").Value = $ws4Vals[$row]
}
